# This script updates the "quiz" marksheet worksheet:
#  - fills in the summary (Right/Wrong/Not-Attempt/Max, Marking, Total) numbers
#    for a graded submission instead of the "Absent" placeholder
#  - fixes the "-1" marking-per-wrong-answer cell so it is stored as a real
#    number instead of text (this is what let float/numeric input break the
#    sheet before)
#  - fills in the student's answers in column A for the answer-key block
#    (rows 16-40), re-using the existing "correct"/"incorrect" cell styles
#  - removes the now-unused extra answer-key columns (D:E for rows 19-40,
#    and the whole third G:H block) so the sheet shrinks back down to A:E

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Summary block (rows 10-12)
# ---------------------------------------------------------------------------

# No. (counts)
$ws.Range("B10").Value() = 18
$ws.Range("C10").Value() = 3
$ws.Range("D10").Value() = 7
$ws.Range("E10").Value() = 28

# Marking (points per right/wrong answer) - keep these numeric, not text
$ws.Range("B11").Value() = 4
$ws.Range("C11").Value() = -1

# Total (score)
$ws.Range("B12").Value() = 72
$ws.Range("C12").Value() = -3
$ws.Range("E12").Value() = "69/112"

# Give the row labels (No./Marking/Total) the same style used by the header
# row above them (A9), instead of the default style.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Answer key block (rows 16-40): fill in the student's answers (column A)
# ---------------------------------------------------------------------------

# Map of row -> student answer ($null means left blank / not attempted)
$studentAnswers = @{
    16 = "Option A"
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    20 = $null
    21 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = $null
    25 = "Option B"
    26 = $null
    27 = "Option A"
    28 = $null
    29 = $null
    30 = "Option B"
    31 = "Option C"
    32 = "Option C"
    33 = "Option D"
    34 = $null
    35 = $null
    36 = "Option A"
    37 = "Option A"
    38 = "Option C"
    39 = "Option D"
    40 = "Option D"
}

foreach ($r in $studentAnswers.Keys) {
    $answer = $studentAnswers[$r]
    $cell = $ws.Range("A$r")
    $correctCell = $ws.Range("B$r")

    if ($answer -eq $null) {
        # leave blank / not attempted - keep the existing "normalStyle" look
        continue
    }

    $cell.Value() = $answer
    $correct = $correctCell.Value()

    if ($answer -eq $correct) {
        # correct answer -> reuse the existing "correctStyle" formatting
        $ws.Range("B10").Copy()
    } else {
        # wrong answer -> reuse the existing "incorrectStyle" formatting
        $ws.Range("C10").Copy()
    }
    $cell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Drop the now-unused extra answer-key columns
# ---------------------------------------------------------------------------

# Third block (columns G:H) is no longer used at all.
$ws.Range("G15:H21").Clear()

# Second block (columns D:E) is only kept for rows 16-18; rows 19-40 drop it.
$ws.Range("D19:E40").Clear()

Write-Host "edit applied"
